$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New shared-string bodies (Cypher queries / filenames / tab name) ----
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE s.clinical_study_designation IN ['UBC01'] and demo.sex in ['Female'] and demo.neutered_indicator in ['No'] 
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$neo4jName = @'
TC34_Canine_StudyUBC-Breed_RespToTreatmt_Sex_NeutStatus_Neo4jData.xlsx
'@
$webName = @'
TC34_Canine_StudyUBC-Breed_RespToTreatmt_Sex_NeutStatus_WebData.xlsx
'@
$studyFilesTab = @'
StudyFilesTab
'@
$caseQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.sex in ['Female'] and demo.neutered_indicator in ['No'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
   coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
'@
$sampleQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and demo.sex in ['Female'] and demo.neutered_indicator in ['No'] 
 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$fileQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and demo.sex in ['Female'] and demo.neutered_indicator in ['No']  
WITH DISTINCT f, parent, c, demo, diag, s
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@
$studyFileQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.sex in ['Female'] and demo.neutered_indicator in ['No']  
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# ---- Row labels (A column) for existing rows are unchanged text ----
$ws.Range("A2").Value = "CasesTab"
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# ---- Populate cells in the exact order the strings were first introduced ----
# Row 2 (Cases): C2 (stat), D2 (neo4j name), E2 (web name) first, then B2 (case query)
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = $neo4jName
$ws.Range("E2").Value = $webName
$ws.Range("A5").Value = $studyFilesTab
$ws.Range("B2").Value = $caseQuery
$ws.Range("B3").Value = $sampleQuery
$ws.Range("B4").Value = $fileQuery
$ws.Range("B5").Value = $studyFileQuery

# ---- Remaining reused cells (stat query / filenames repeated on each row) ----
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $neo4jName
$ws.Range("E3").Value = $webName
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = $neo4jName
$ws.Range("E4").Value = $webName
$ws.Range("C5").Value = $statQuery
$ws.Range("D5").Value = $neo4jName
$ws.Range("E5").Value = $webName

# ---- Fonts / wrap styling ----
# Columns B:C (the long query text) cells: font size 13 + wrap text (style created first)
$ws.Range("B2:C5").Font.Size = 13
$ws.Range("B2:C5").WrapText = $true

# Column A (tab-name) cells for the data rows: font size 13, no wrap (style created second)
$ws.Range("A2:A5").Font.Size = 13

# ---- Row heights ----
$ws.Rows(2).RowHeight = 396.75
$ws.Rows(3).RowHeight = 362.25
$ws.Rows(4).RowHeight = 409.5
$ws.Rows(5).RowHeight = 409.5

# ---- Sheet view: selection moves to C2 ----
$ws.Range("C2").Select() | Out-Null
